$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.303.35"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.94"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.23"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4723"
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2869"
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06471"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.63"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07779"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "95.97"
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.866.60"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7149"
$ws.Range("E14").Value = "  -3.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.115"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.41"
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.290.44"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.95"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.110.76"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.235"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.35"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.964"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.875"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09592"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.312"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.480"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.197"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.110"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04797"
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.115"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("E36").Value = "  -2.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01875"
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("E39").Value = "  +2.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.20"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.205"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("E42").Value = "  -4.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4184"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8219"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.43"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.586"
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.984"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.97"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "891.62"
$ws.Range("E50").Value = "  -2.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05740"
$ws.Range("E51").Value = "  +0.09%  "
